$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fn1"
$ws.Range("C2").Value = "Tshr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.116717
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.960636
$ws.Range("N2").Value = 2.881908
$ws.Range("O2").Value = 0.2124324572954377
$ws.Range("P2").Value = 0.2124324572954377
$ws.Range("Q2").Value = 25.974546184004
$ws.Range("R2").Value = 233.770915656036
$ws.Range("S2").Value = 0.0150746070080027
$ws.Range("T2").Value = 0.01507460700800271

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fn1"
$ws.Range("C3").Value = "Tshr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.116717
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.5285266666666667
$ws.Range("N3").Value = 1.58558
$ws.Range("O3").Value = 0.116876963330717
$ws.Range("P3").Value = 0.116876963330717
$ws.Range("Q3").Value = 14.29078268231778
$ws.Range("R3").Value = 128.61704414086
$ws.Range("S3").Value = 0.008293809302638712
$ws.Range("T3").Value = 0.008293809302638714

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fn1"
$ws.Range("C4").Value = "Tshr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.116717
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.032914666666667
$ws.Range("N4").Value = 9.098744
$ws.Range("O4").Value = 0.6706905793738452
$ws.Range("P4").Value = 0.6706905793738454
$ws.Range("Q4").Value = 82.00669356704977
$ws.Range("R4").Value = 738.060242103448
$ws.Range("S4").Value = 0.04759346587969586
$ws.Range("T4").Value = 0.04759346587969587

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fn1"
$ws.Range("C5").Value = "Tshr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.960636
$ws.Range("N5").Value = 2.881908
$ws.Range("O5").Value = 0.2124324572954377
$ws.Range("P5").Value = 0.2124324572954377
$ws.Range("Q5").Value = 331.963696184244
$ws.Range("R5").Value = 2987.673265658196
$ws.Range("S5").Value = 0.1926586984600814
$ws.Range("T5").Value = 0.1926586984600815

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fn1"
$ws.Range("C6").Value = "Tshr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.5285266666666667
$ws.Range("N6").Value = 1.58558
$ws.Range("O6").Value = 0.116876963330717
$ws.Range("P6").Value = 0.116876963330717
$ws.Range("Q6").Value = 182.6411521102733
$ws.Range("R6").Value = 1643.77036899246
$ws.Range("S6").Value = 0.1059977553427576
$ws.Range("T6").Value = 0.1059977553427576

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fn1"
$ws.Range("C7").Value = "Tshr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.032914666666667
$ws.Range("N7").Value = 9.098744
$ws.Range("O7").Value = 0.6706905793738452
$ws.Range("P7").Value = 0.6706905793738454
$ws.Range("Q7").Value = 1048.073945758925
$ws.Range("R7").Value = 9432.665511830328
$ws.Range("S7").Value = 0.6082609773321962
$ws.Range("T7").Value = 0.6082609773321964

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fn1"
$ws.Range("C8").Value = "Tshr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.960636
$ws.Range("N8").Value = 2.881908
$ws.Range("O8").Value = 0.2124324572954377
$ws.Range("P8").Value = 0.2124324572954377
$ws.Range("Q8").Value = 8.096949797792002
$ws.Range("R8").Value = 72.872548180128
$ws.Range("S8").Value = 0.004699151827353542
$ws.Range("T8").Value = 0.004699151827353545

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fn1"
$ws.Range("C9").Value = "Tshr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.5285266666666667
$ws.Range("N9").Value = 1.58558
$ws.Range("O9").Value = 0.116876963330717
$ws.Range("P9").Value = 0.116876963330717
$ws.Range("Q9").Value = 4.454813151697778
$ws.Range("R9").Value = 40.09331836528001
$ws.Range("S9").Value = 0.002585398685320708
$ws.Range("T9").Value = 0.002585398685320708

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fn1"
$ws.Range("C10").Value = "Tshr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.032914666666667
$ws.Range("N10").Value = 9.098744
$ws.Range("O10").Value = 0.6706905793738452
$ws.Range("P10").Value = 0.6706905793738454
$ws.Range("Q10").Value = 25.56364512363378
$ws.Range("R10").Value = 230.072806112704
$ws.Range("S10").Value = 0.01483613616195315
$ws.Range("T10").Value = 0.01483613616195316

